# The sheet holds yearly rows starting at row 2 (2008, 2009, 2010, 2011, ...).
# The target keeps only the 2010/2011 rows (shifted up to rows 2/3), i.e. the
# 2008 and 2009 rows are removed and everything below shifts up accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A3").EntireRow.Delete()
